$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the data range to Text format before writing, so numeric-looking
# strings (e.g. "250.64") are kept as literal text instead of being
# auto-converted to floating point numbers by the Value setter.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "37.387.92"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.095.89"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "250.64"
$ws.Range("E5").Value = "  +1.58%  "
$ws.Range("D6").Value = "0.659"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "52.11"
$ws.Range("E8").Value = "  +14.49%  "
$ws.Range("D9").Value = "61.96"
$ws.Range("E9").Value = "  +9.75%  "
$ws.Range("E10").Value = "  +3.23%  "
$ws.Range("E11").Value = "  +3.94%  "
$ws.Range("D12").Value = "0.106"
$ws.Range("E12").Value = "  +6.80%  "
$ws.Range("D13").Value = "15.16"
$ws.Range("E13").Value = "  +2.87%  "
$ws.Range("D14").Value = "2.401.40"
$ws.Range("E14").Value = "  +4.19%  "
$ws.Range("D15").Value = "0.831"
$ws.Range("E15").Value = "  +3.33%  "
$ws.Range("D16").Value = "2.099.71"
$ws.Range("E16").Value = "  +3.69%  "
$ws.Range("D17").Value = "5.11"
$ws.Range("E17").Value = "  +4.59%  "
$ws.Range("D18").Value = "37.292.08"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "72.04"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "13.91"
$ws.Range("E20").Value = "  +7.52%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +2.72%  "
$ws.Range("D22").Value = "240.42"
$ws.Range("E22").Value = "  +2.84%  "
$ws.Range("D23").Value = "5.21"
$ws.Range("E23").Value = "  +4.51%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "2.46"
$ws.Range("E25").Value = "  -0.76%  "
$ws.Range("D26").Value = "170.50"
$ws.Range("E26").Value = "  +5.34%  "
$ws.Range("D27").Value = "9.19"
$ws.Range("E27").Value = "  +8.03%  "
$ws.Range("D28").Value = "20.73"
$ws.Range("E28").Value = "  +5.13%  "
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("D30").Value = "0.123"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "1.07"
$ws.Range("E31").Value = "  +26.19%  "
$ws.Range("E32").Value = "  +2.45%  "
$ws.Range("D33").Value = "0.0610"
$ws.Range("E33").Value = "  +4.40%  "
$ws.Range("E34").Value = "  +8.66%  "
$ws.Range("D35").Value = "19.85"
$ws.Range("E35").Value = "  -6.52%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "2.31"
$ws.Range("E37").Value = "  +4.96%  "
$ws.Range("D38").Value = "1.85"
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("E40").Value = "  -2.89%  "
$ws.Range("D41").Value = "18.36"
$ws.Range("E41").Value = "  +12.22%  "
$ws.Range("E42").Value = "  +3.97%  "
$ws.Range("E43").Value = "  +7.04%  "
$ws.Range("D44").Value = "99.35"
$ws.Range("E44").Value = "  +2.17%  "
$ws.Range("B45").Value = "HuobiToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D45").Value = "2.75"
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.0902"
$ws.Range("E46").Value = "  +11.16%  "
$ws.Range("D47").Value = "3.01"
$ws.Range("E47").Value = "  +8.63%  "
$ws.Range("D48").Value = "1.321.67"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "7.02"
$ws.Range("E49").Value = "  +15.20%  "
$ws.Range("D50").Value = "2.288.49"
$ws.Range("E50").Value = "  +3.82%  "
$ws.Range("E51").Value = "  +1.62%  "

# Restore original (default/general) formatting now that the literal text
# values are committed, so cell styling matches the un-formatted source.
$dataRange.ClearFormats()
